$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-05-15"

# Update the header label for the "through" column
$ws.Range("I1").Value = "2022 (through 05-15)"

# Update the May data point and the Total row for the 2022 column
$ws.Range("I6").Value = 57
$ws.Range("I14").Value = 609
